# Auto-generated edit script: updates Price (D) and Volume(1h) (E)
# columns for rows 2-51 on Sheet1, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.432.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.952.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.75%  '

$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.44'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.92%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.72'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.93%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0855'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.78%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.237.57'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -11.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.821'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.97%  '

$ws.Range("E16").Value = '  -4.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.951.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.357.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0887'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.60%  '

$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("E24").Value = '  -7.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("E28").Value = '  +2.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.46%  '

$ws.Range("E30").Value = '  -1.74%  '

$ws.Range("E31").Value = '  -1.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0650'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.30'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.65%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.15'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0983'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.29%  '

$ws.Range("E41").Value = '  +0.46%  '

$ws.Range("E42").Value = '  -6.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0212'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.356.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.09%  '

$ws.Range("E46").Value = '  -6.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.91%  '

$ws.Range("E49").Value = '  -0.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.20'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.128.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.05%  '
